# Auto-generated data-driven update of Hades_Profits market-data sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{ Row=80; H=384.75674; I=284.15384; J=622.5455; K=852.4615200000001; L=1867.6365; M=145.5384799999999; N=-3863.6365 }
    @{ Row=83; H=384.75674; I=284.15384; J=622.5455; K=2557.38456; L=5602.9095; M=2434.61544; N=-15586.9095 }
    @{ Row=88; H=21200.3; I=602; J=26349.875; K=602; L=26349.875; M=-196; N=-27161.875 }
    @{ Row=91; H=21200.3; I=602; J=26349.875; K=602; L=26349.875; M=802; N=-29157.875 }
    @{ Row=93; H=34967; I=0; J=34967; K=0; L=34967; N=-39959 }
    @{ Row=116; H=2071.4285; I=1700; J=2350; K=1700; L=2350; M=1742; N=-9234 }
    @{ Row=129; H=905.74; I=540.3; J=997.1; K=1620.9; L=2991.3; M=3379.1; N=-12991.3 }
    @{ Row=132; H=1486551.5; I=1730.6552; J=12251502; K=5191.9656; L=36754506; M=-2661.9656; N=-36759566 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{ Row=61; H=83502210; I=166834720; J=169696.67; K=166834720; L=169696.67; M=-166834508; N=-170120.67 }
    @{ Row=88; H=6650.273; I=3833.875; J=8259.643; K=3833.875; L=8259.643; M=-3427.875; N=-9071.643 }
    @{ Row=91; H=6650.273; I=3833.875; J=8259.643; K=3833.875; L=8259.643; M=-2429.875; N=-11067.643 }
    @{ Row=132; H=51115.883; I=31095.605; J=124523.555; K=93286.815; L=373570.665; M=-90756.815; N=-378630.665 }
    @{ Row=136; H=83502210; I=166834720; J=169696.67; K=500504160; L=509090.01; M=-500501610; N=-514190.01 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{ Row=107; H=1264.15; I=1293.0588; J=1100.3334; K=1293.0588; L=1100.3334; M=626.9412; N=-4940.3334 }
    @{ Row=134; H=3393.2285; I=2652.64; J=5244.7; K=7957.92; L=15734.1; M=-5422.92; N=-20804.1 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{ Row=31; H=3358.1072; I=1712.6111; J=6320; K=1712.6111; L=6320; M=-1417.6111; N=-6910 }
    @{ Row=34; H=3358.1072; I=1712.6111; J=6320; K=1712.6111; L=6320; M=-1510.6111; N=-6724 }
    @{ Row=75; H=38000; I=0; J=38000; K=0; L=38000; N=-39996 }
    @{ Row=78; H=38000; I=0; J=38000; K=0; L=114000; N=-123984 }
    @{ Row=94; H=4616.154; I=12153; J=1266.4445; K=12153; L=1266.4445; M=-11702; N=-2168.4445 }
    @{ Row=99; H=2623.72; I=2208.8; J=2900.3333; K=2208.8; L=2900.3333; M=-710.8000000000002; N=-5896.3333 }
    @{ Row=126; H=2623.72; I=2208.8; J=2900.3333; K=6626.400000000001; L=8700.999899999999; M=-4156.400000000001; N=-13640.9999 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{ Row=19; H=0; I=0; J=0; K=0; L=0; N=$null }
    @{ Row=76; H=3462.5; I=1400; J=3757.1428; K=4200; L=11271.4284; M=-3817; N=-12037.4284 }
    @{ Row=79; H=3462.5; I=1400; J=3757.1428; K=4200; L=11271.4284; M=-2874; N=-13923.4284 }
    @{ Row=131; H=980.9259; I=514.2222; J=1039.2639; K=1542.6666; L=3117.7917; M=3497.3334; N=-13197.7917 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{ Row=132; H=85174.086; I=72489; J=102933.2; K=217467; L=308799.6; M=-214937; N=-313859.6 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{ Row=136; H=129750.19; I=112700.336; J=151671.42; K=338101.008; L=455014.26; M=-335551.008; N=-460114.26 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{ Row=46; H=30429; I=0; J=30429; K=0; L=30429; N=-30891 }
    @{ Row=74; H=0; I=0; J=0; K=0; L=0; N=$null }
    @{ Row=77; H=0; I=0; J=0; K=0; L=0; N=$null }
    @{ Row=113; H=803.34375; I=490.72726; J=1491.1; K=1472.18178; L=4473.299999999999; M=697.8182200000001; N=-8813.3 }
    @{ Row=132; H=50243.023; I=38080.223; J=73699.86; K=114240.669; L=221099.58; M=-111710.669; N=-226159.58 }
    @{ Row=134; H=30429; I=0; J=30429; K=0; L=91287; N=-96357 }
    @{ Row=136; H=49819.023; I=26553.334; J=503500; K=79660.002; L=1510500; M=-77110.002; N=-1515600 }
)
foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("H","I","J","K","L","M","N")) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $val = $u[$col]
            if ($null -eq $val) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}
